$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.632.59"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.96%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.879.11"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.95%  "

$ws.Range("E4").Value = "  +0.34%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.07"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.31%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.008"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.45%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5107"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.52%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3931"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.17%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08416"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.89%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.113"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.35%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.77"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.35%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.279"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.75%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.880.91"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.98%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.51"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.281"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.47%  "

$ws.Range("E16").Value = "  +0.21%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001107"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.37%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.54"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.14%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06732"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.76"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.14%  "

$ws.Range("E21").Value = "  +0.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.976"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.641.80"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.81%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.16"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.30%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.246"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.23%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.097.88"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.70"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.72%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.83"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.376"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.99"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.27%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1057"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.53%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.058"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.33%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.823"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.622"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.89%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02461"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06551"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.76%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2190"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.57%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.917"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.60%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.271"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.96%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.199"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.31%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6486"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.31%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.084"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.22%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.21"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.007"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.41%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6079"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.86%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.02"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.11%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.698"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.10%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.037"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.87%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.221"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.73%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "122.52"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.37%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.189"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -6.65%  "
